$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is a new compound entry ("Rhodamine_n"), mirroring row 4's data/format
# except for Compound name (A), pKa (L) and chemcharge (N).

$ws.Range("A5").NumberFormat = "0.00"
$ws.Range("A5").Value = "Rhodamine_n"

$ws.Range("B5").NumberFormat = "0.00"
$ws.Range("B5").Value = 566.99

$ws.Range("C5").Value = [double]"7.8946170041811414E-6"

$ws.Range("D5").NumberFormat = "0.00E+00"
$ws.Range("D5").Value = [double]"1.2887999999999998E-2"

$ws.Range("E5").NumberFormat = "0.00E+00"
$ws.Range("E5").Value = [double]"4.9712889999999998E-10"

$ws.Range("F5").NumberFormat = "0.00"
$ws.Range("F5").Value = 1440

$ws.Range("G5").NumberFormat = "0.00"
$ws.Range("G5").Value = 2880

$ws.Range("H5").NumberFormat = "0.00"
$ws.Range("H5").Value = 144

$ws.Range("I5").NumberFormat = "0.00"
$ws.Range("I5").Value = -1.33

$ws.Range("J5").NumberFormat = "0.00"
$ws.Range("J5").Value = -33.538476999866809

$ws.Range("K5").NumberFormat = "0.00"
$ws.Range("K5").Value = 4.3643633546157306

$ws.Range("L5").NumberFormat = "0.00"
$ws.Range("L5").Value = 999

$ws.Range("M5").NumberFormat = "0.00"

$ws.Range("N5").NumberFormat = "0.00"
$ws.Range("N5").Value = 0

$ws.Range("O5").NumberFormat = "0.00"
$ws.Range("O5").Value = 16.437000000000001

$ws.Range("P5").NumberFormat = "0.00"
$ws.Range("P5").Value = 2.33

$ws.Range("Q5").NumberFormat = "0.00"
$ws.Range("Q5").Value = 0

$ws.Range("R5").NumberFormat = "0.00"
$ws.Range("R5").Value = 1.47

$ws.Range("S5").NumberFormat = "0.00"
$ws.Range("S5").Value = 3.7159

$ws.Range("T5").NumberFormat = "0.00"
$ws.Range("T5").Value = "CCN(CC)C1=CC2=C(C=C1)C(=C3C=CC(=[N+](CC)CC)C=C3O2)C4=C(C=C(C=C4)C(=O)[O-])C(=O)[O-].[Na+].[Na+].[Cl-]"

$ws.Range("U5").Value = "37299-86-8"

# Update active selection to match the author's final cursor position
$ws.Range("N10").Select()
